$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 95) {
        # Row 95 is a special case: remaining resets to 10 and the start date moves forward
        $ws.Cells.Item($row, 5).Value = 10
        $ws.Cells.Item($row, 6).Value = 20260301
        continue
    }

    if ($row -eq 36) {
        continue
    }
    $current = $ws.Cells.Item($row, 5).Value2
    if ($current -eq $null) {
        continue
    }
    $ws.Cells.Item($row, 5).Value = $current - 1
}
